# Append: 2026-01-18 12:36 JST
# Target sheet is the "ランサーズ" sheet (first worksheet / rId1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2026-01-18 12:36:14"

# --- Drop every existing hyperlink up front -------------------------------
# (the engine's Range.Hyperlinks.Delete() clears the whole sheet's
# hyperlink collection, which is exactly what we want here since we are
# going to rebuild all four of them, in row order, below).
$ws.Range("A1").Hyperlinks.Delete()

# --- Row 2: just refresh the "fetched at" timestamp -----------------------
$ws.Range("A2").Value = $newTimestamp

# --- Shift the existing row 3 ("製造業DX...") down to row 4 ---------------
# Insert() pushes row 3 (and its formatting) down to row 4, leaving a blank
# row 3 behind for the brand-new listing.
$ws.Rows.Item(3).Insert()

# refresh timestamp on the row that just moved down to row 4
$ws.Range("A4").Value = $newTimestamp

# --- New row 3: 日本Amazonアカウント SP-APIの取得および権限付与と設定 -----
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "日本Amazonアカウント SP-APIの取得および権限付与と設定"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5473858"
$ws.Range("G3").Value = 180
$ws.Range("H3").Value = "🔥API"

# --- New row 5: 【急募】プログラム修正依頼!スキルを活かしてみませんか? ---
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5473840"
$ws.Range("G5").Value = 13

# --- Rebuild the hyperlinks, in row order, so relationship ids line up ----
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473648")
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5473858")
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5468432")
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5473840")
$ws.Range("F5").Style = "Hyperlink"

Write-Output "applied 2026-01-18 12:36 append"
